$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row after "Contact" (row 10) for "Jurisdiction"
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
